# Append new codebook entries (sd_imputed, yi, vi, control_accuracy) to Sheet1,
# reflecting the updated models (quadratic effect of accuracy) described in the
# commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new variable names (column A) first, top to bottom.
$ws.Range("A77").Value = "sd_imputed"
$ws.Range("A78").Value = "yi"
$ws.Range("A79").Value = "vi"
$ws.Range("A80").Value = "control_accuracy"

# Then the descriptions, starting with control_accuracy (row 80) ...
$ws.Range("B80").Value = "Accuracy in the control condition. Because this variable is calculated in the wrangling code, it is not included in the raw or cleaned data. When accuracy is reported as a proportion, this is simply the mean accuracy rate in the control condition. When accuracy is reported as a count, this value is estimated by dividing the reported mean count in the control condition by the number of critical items on the test."

# ... then sd_imputed's description and coded values (row 77) ...
$ws.Range("B77").Value = "Indicator for whether the SD for accuracy has been imputed."
$ws.Range("D77").Value = "0 = not imputed; 1 = imputed"

# ... then yi (row 78) ...
$ws.Range("B78").Value = "Standardized mean difference for the misinformation effect."

# ... then vi (row 79).
$ws.Range("B79").Value = "Sampling variance for the standardized mean difference."

# Fill in the "type" column (numeric) for all four new rows.
$ws.Range("C77").Value = "numeric"
$ws.Range("C78").Value = "numeric"
$ws.Range("C79").Value = "numeric"
$ws.Range("C80").Value = "numeric"

$ws.Range("B78").Select() | Out-Null
